$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sites")

# New header cells (bold, same style as existing headers in row 1)
$ws.Range("D1").Value = "Study Phase"
$ws.Range("E1").Value = "Status"

# New data cells in row 2 (normal style)
$ws.Range("D2").Value = "Phase II/III"
$ws.Range("E2").Value = "Planning"
